# Generate Report for Handoff
# Updates the localization-status report:
#  - Priority changes from "low" to "ht" for the last four rows (rows 4-7)
#    on both the zh-cn and de-de sheets.
#  - Latest Handoff Datetime on the zh-cn sheet (rows 4-7) is refreshed.
#  - Latest Handoff Datetime on the de-de sheet (rows 4-7), together with the
#    matching "Latest HO Xliff Generate Date" on the Overview sheet (which
#    shares the same underlying value), is refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-26 12:30:59"

    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-26 12:31:11"

    $wsOverview.Range("G$r").Value = "2016-08-26 12:31:11"
}
